$wb = $excel.ActiveWorkbook

# Rename the "KNX Group Addresses" sheet to "KNX GAs".
$ws = $wb.Worksheets.Item("KNX Group Addresses")
$ws.Name = "KNX GAs"
